# ULYSSES-6858: Better fix for whitespace that doesn't break other styles
#
# Gives the separator <w:tab/> run inside the footnote / endnote reference
# paragraph the same "footnote reference" / "endnote reference" character
# style (plus matching font/size) that the reference mark run already has,
# instead of leaving it with no formatting at all.

$d = $word.ActiveDocument

# --- Footnote --------------------------------------------------------
$fn = $d.Footnotes(1)
$tabRun = $fn.Range.Duplicate
$tabRun.SetRange(0, 1)
$tabRun.Font.NameAscii = "Helvetica"
$tabRun.Font.NameFarEast = "Helvetica"
$tabRun.Font.NameOther = "Helvetica"
$tabRun.Font.Size = 12
$tabRun.Style = "footnote reference"

# --- Endnote -----------------------------------------------------------
$en = $d.Endnotes(1)
$tabRun2 = $en.Range.Duplicate
$tabRun2.SetRange(0, 1)
$tabRun2.Font.NameAscii = "Helvetica"
$tabRun2.Font.NameFarEast = "Helvetica"
$tabRun2.Font.NameOther = "Helvetica"
$tabRun2.Font.Size = 12
$tabRun2.Style = "endnote reference"
